# Actualizacion README - Indicacion Ruta de Script
#
# Mark rows 6-12 of the "STATUS" column (Q) as "ATENDIDO" instead of
# "PENDIENTE", and move the sheet's active selection to P13 (reflecting
# where the user last left the cursor after reviewing the updated rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the STATUS column (Q) for rows 6 through 12.
for ($r = 6; $r -le 12; $r++) {
    $ws.Cells.Item($r, 17).Value = "ATENDIDO"
}

# Leave the cursor/selection on P13, matching where the review ended.
$ws.Range("P13").Select()
